$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 and 5 each describe one bat observation. The two records get
# swapped: what used to be row 5's data moves to row 4, and vice versa.

$ws.Range("A4").Value = 131117036
$ws.Range("B4").Value = 56748
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 205998
$ws.Range("F4").Value = "Nordfladdermus"
$ws.Range("G4").Value = "Eptesicus nilssonii"
$ws.Range("H4").Value = "(A.Keyserling & Blasius, 1839)"

# I4/I5 hold digit-only text ("42"/"443"), not numbers, in the source file.
# Force text storage so Excel doesn't silently re-type them as numeric.
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "443"
$ws.Range("I4").Style = "Normal"

$ws.Range("A5").Value = 131116964
$ws.Range("B5").Value = 56762
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 100092
$ws.Range("F5").Value = "Större brunfladdermus"
$ws.Range("G5").Value = "Nyctalus noctula"
$ws.Range("H5").Value = "(Schreber, 1774)"

$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "42"
$ws.Range("I5").Style = "Normal"
